$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three rows (old rows 53-55); shrinks table to A1:E52
$ws.Rows("53:55").Delete()

# Update auto-computed columns: B (time in ms), C (auto capacity), D (auto scs keys), E (probabilities)
$ws.Cells.Item(2, 2).Value = 33.94444444444444
$ws.Cells.Item(2, 3).Value = 1.95
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0.148
$ws.Cells.Item(3, 2).Value = 33.94444444444444
$ws.Cells.Item(3, 3).Value = 1.95
$ws.Cells.Item(3, 4).Value = 2
$ws.Cells.Item(3, 5).Value = 0.003
$ws.Cells.Item(4, 2).Value = 33.94444444444444
$ws.Cells.Item(4, 3).Value = 1.95
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0.005
$ws.Cells.Item(5, 2).Value = 33.94444444444444
$ws.Cells.Item(5, 3).Value = 1.95
$ws.Cells.Item(5, 4).Value = 4
$ws.Cells.Item(5, 5).Value = 0.009000000000000001
$ws.Cells.Item(6, 2).Value = 33.94444444444444
$ws.Cells.Item(6, 3).Value = 1.95
$ws.Cells.Item(6, 4).Value = 5
$ws.Cells.Item(6, 5).Value = 0.021
$ws.Cells.Item(7, 2).Value = 33.94444444444444
$ws.Cells.Item(7, 3).Value = 1.95
$ws.Cells.Item(7, 4).Value = 6
$ws.Cells.Item(7, 5).Value = 0.032
$ws.Cells.Item(8, 2).Value = 33.94444444444444
$ws.Cells.Item(8, 3).Value = 1.95
$ws.Cells.Item(8, 4).Value = 7
$ws.Cells.Item(8, 5).Value = 0.049
$ws.Cells.Item(9, 2).Value = 33.94444444444444
$ws.Cells.Item(9, 3).Value = 1.95
$ws.Cells.Item(9, 4).Value = 8
$ws.Cells.Item(9, 5).Value = 0.057
$ws.Cells.Item(10, 2).Value = 33.94444444444444
$ws.Cells.Item(10, 3).Value = 1.95
$ws.Cells.Item(10, 4).Value = 9
$ws.Cells.Item(10, 5).Value = 0.052
$ws.Cells.Item(11, 2).Value = 33.94444444444444
$ws.Cells.Item(11, 3).Value = 1.95
$ws.Cells.Item(11, 4).Value = 10
$ws.Cells.Item(11, 5).Value = 0.038
$ws.Cells.Item(12, 2).Value = 33.94444444444444
$ws.Cells.Item(12, 3).Value = 1.95
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = 0.029
$ws.Cells.Item(13, 2).Value = 33.94444444444444
$ws.Cells.Item(13, 3).Value = 1.95
$ws.Cells.Item(13, 4).Value = 12
$ws.Cells.Item(13, 5).Value = 0.036
$ws.Cells.Item(14, 2).Value = 33.94444444444444
$ws.Cells.Item(14, 3).Value = 1.95
$ws.Cells.Item(14, 4).Value = 13
$ws.Cells.Item(14, 5).Value = 0.028
$ws.Cells.Item(15, 2).Value = 33.94444444444444
$ws.Cells.Item(15, 3).Value = 1.95
$ws.Cells.Item(15, 4).Value = 14
$ws.Cells.Item(15, 5).Value = 0.021
$ws.Cells.Item(16, 2).Value = 33.94444444444444
$ws.Cells.Item(16, 3).Value = 1.95
$ws.Cells.Item(16, 4).Value = 15
$ws.Cells.Item(16, 5).Value = 0.025
$ws.Cells.Item(17, 2).Value = 33.94444444444444
$ws.Cells.Item(17, 3).Value = 1.95
$ws.Cells.Item(17, 4).Value = 16
$ws.Cells.Item(17, 5).Value = 0.037
$ws.Cells.Item(18, 2).Value = 33.94444444444444
$ws.Cells.Item(18, 3).Value = 1.95
$ws.Cells.Item(18, 4).Value = 17
$ws.Cells.Item(18, 5).Value = 0.039
$ws.Cells.Item(19, 2).Value = 33.94444444444444
$ws.Cells.Item(19, 3).Value = 1.95
$ws.Cells.Item(19, 4).Value = 18
$ws.Cells.Item(19, 5).Value = 0.031
$ws.Cells.Item(20, 2).Value = 33.94444444444444
$ws.Cells.Item(20, 3).Value = 1.95
$ws.Cells.Item(20, 4).Value = 19
$ws.Cells.Item(20, 5).Value = 0.043
$ws.Cells.Item(21, 2).Value = 33.94444444444444
$ws.Cells.Item(21, 3).Value = 1.95
$ws.Cells.Item(21, 4).Value = 20
$ws.Cells.Item(21, 5).Value = 0.022
$ws.Cells.Item(22, 2).Value = 33.94444444444444
$ws.Cells.Item(22, 3).Value = 1.95
$ws.Cells.Item(22, 4).Value = 21
$ws.Cells.Item(22, 5).Value = 0.025
$ws.Cells.Item(23, 2).Value = 33.94444444444444
$ws.Cells.Item(23, 3).Value = 1.95
$ws.Cells.Item(23, 4).Value = 22
$ws.Cells.Item(23, 5).Value = 0.02
$ws.Cells.Item(24, 2).Value = 33.94444444444444
$ws.Cells.Item(24, 3).Value = 1.95
$ws.Cells.Item(24, 4).Value = 23
$ws.Cells.Item(24, 5).Value = 0.024
$ws.Cells.Item(25, 2).Value = 33.94444444444444
$ws.Cells.Item(25, 3).Value = 1.95
$ws.Cells.Item(25, 4).Value = 24
$ws.Cells.Item(25, 5).Value = 0.018
$ws.Cells.Item(26, 2).Value = 33.94444444444444
$ws.Cells.Item(26, 3).Value = 1.95
$ws.Cells.Item(26, 4).Value = 25
$ws.Cells.Item(26, 5).Value = 0.018
$ws.Cells.Item(27, 2).Value = 33.94444444444444
$ws.Cells.Item(27, 3).Value = 1.95
$ws.Cells.Item(27, 4).Value = 26
$ws.Cells.Item(27, 5).Value = 0.011
$ws.Cells.Item(28, 2).Value = 33.94444444444444
$ws.Cells.Item(28, 3).Value = 1.95
$ws.Cells.Item(28, 4).Value = 27
$ws.Cells.Item(28, 5).Value = 0.019
$ws.Cells.Item(29, 2).Value = 33.94444444444444
$ws.Cells.Item(29, 3).Value = 1.95
$ws.Cells.Item(29, 4).Value = 28
$ws.Cells.Item(29, 5).Value = 0.014
$ws.Cells.Item(30, 2).Value = 33.94444444444444
$ws.Cells.Item(30, 3).Value = 1.95
$ws.Cells.Item(30, 4).Value = 29
$ws.Cells.Item(30, 5).Value = 0.012
$ws.Cells.Item(31, 2).Value = 33.94444444444444
$ws.Cells.Item(31, 3).Value = 1.95
$ws.Cells.Item(31, 4).Value = 30
$ws.Cells.Item(31, 5).Value = 0.011
$ws.Cells.Item(32, 2).Value = 33.94444444444444
$ws.Cells.Item(32, 3).Value = 1.95
$ws.Cells.Item(32, 4).Value = 31
$ws.Cells.Item(32, 5).Value = 0.011
$ws.Cells.Item(33, 2).Value = 33.94444444444444
$ws.Cells.Item(33, 3).Value = 1.95
$ws.Cells.Item(33, 4).Value = 32
$ws.Cells.Item(33, 5).Value = 0.009000000000000001
$ws.Cells.Item(34, 2).Value = 33.94444444444444
$ws.Cells.Item(34, 3).Value = 1.95
$ws.Cells.Item(34, 4).Value = 33
$ws.Cells.Item(34, 5).Value = 0.011
$ws.Cells.Item(35, 2).Value = 33.94444444444444
$ws.Cells.Item(35, 3).Value = 1.95
$ws.Cells.Item(35, 4).Value = 34
$ws.Cells.Item(35, 5).Value = 0.007
$ws.Cells.Item(36, 2).Value = 33.94444444444444
$ws.Cells.Item(36, 3).Value = 1.95
$ws.Cells.Item(36, 4).Value = 35
$ws.Cells.Item(36, 5).Value = 0.008
$ws.Cells.Item(37, 2).Value = 33.94444444444444
$ws.Cells.Item(37, 3).Value = 1.95
$ws.Cells.Item(37, 4).Value = 36
$ws.Cells.Item(37, 5).Value = 0.012
$ws.Cells.Item(38, 2).Value = 33.94444444444444
$ws.Cells.Item(38, 3).Value = 1.95
$ws.Cells.Item(38, 4).Value = 37
$ws.Cells.Item(38, 5).Value = 0.006
$ws.Cells.Item(39, 2).Value = 33.94444444444444
$ws.Cells.Item(39, 3).Value = 1.95
$ws.Cells.Item(39, 4).Value = 38
$ws.Cells.Item(39, 5).Value = 0.004
$ws.Cells.Item(40, 2).Value = 33.94444444444444
$ws.Cells.Item(40, 3).Value = 1.95
$ws.Cells.Item(40, 4).Value = 39
$ws.Cells.Item(40, 5).Value = 0.006
$ws.Cells.Item(41, 2).Value = 33.94444444444444
$ws.Cells.Item(41, 3).Value = 1.95
$ws.Cells.Item(41, 4).Value = 40
$ws.Cells.Item(41, 5).Value = 0.002
$ws.Cells.Item(42, 2).Value = 33.94444444444444
$ws.Cells.Item(42, 3).Value = 1.95
$ws.Cells.Item(42, 4).Value = 41
$ws.Cells.Item(42, 5).Value = 0.001
$ws.Cells.Item(43, 2).Value = 33.94444444444444
$ws.Cells.Item(43, 3).Value = 1.95
$ws.Cells.Item(43, 4).Value = 42
$ws.Cells.Item(43, 5).Value = 0.001
$ws.Cells.Item(44, 2).Value = 33.94444444444444
$ws.Cells.Item(44, 3).Value = 1.95
$ws.Cells.Item(44, 4).Value = 43
$ws.Cells.Item(44, 5).Value = 0.003
$ws.Cells.Item(45, 2).Value = 33.94444444444444
$ws.Cells.Item(45, 3).Value = 1.95
$ws.Cells.Item(45, 4).Value = 44
$ws.Cells.Item(45, 5).Value = 0.005
$ws.Cells.Item(46, 2).Value = 33.94444444444444
$ws.Cells.Item(46, 3).Value = 1.95
$ws.Cells.Item(46, 4).Value = 46
$ws.Cells.Item(46, 5).Value = 0.003
$ws.Cells.Item(47, 2).Value = 33.94444444444444
$ws.Cells.Item(47, 3).Value = 1.95
$ws.Cells.Item(47, 4).Value = 47
$ws.Cells.Item(47, 5).Value = 0.005
$ws.Cells.Item(48, 2).Value = 33.94444444444444
$ws.Cells.Item(48, 3).Value = 1.95
$ws.Cells.Item(48, 4).Value = 49
$ws.Cells.Item(48, 5).Value = 0.003
$ws.Cells.Item(49, 2).Value = 33.94444444444444
$ws.Cells.Item(49, 3).Value = 1.95
$ws.Cells.Item(49, 4).Value = 50
$ws.Cells.Item(49, 5).Value = 0.001
$ws.Cells.Item(50, 2).Value = 33.94444444444444
$ws.Cells.Item(50, 3).Value = 1.95
$ws.Cells.Item(50, 4).Value = 52
$ws.Cells.Item(50, 5).Value = 0.001
$ws.Cells.Item(51, 2).Value = 33.94444444444444
$ws.Cells.Item(51, 3).Value = 1.95
$ws.Cells.Item(51, 4).Value = 55
$ws.Cells.Item(51, 5).Value = 0.002
$ws.Cells.Item(52, 2).Value = 33.94444444444444
$ws.Cells.Item(52, 3).Value = 1.95
$ws.Cells.Item(52, 4).Value = 59
$ws.Cells.Item(52, 5).Value = 0.001
